$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price/volume text values are written as text (matching original inlineStr cells)
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "58.953.81"
$ws.Range("E2").Value = "  +0.17%  "
$ws.Range("D3").Value = "2.506.73"
$ws.Range("E3").Value = "  +0.47%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "533.15"
$ws.Range("E5").Value = "  -0.28%  "
$ws.Range("D6").Value = "135.72"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  +0.30%  "
$ws.Range("E8").Value = "  +0.45%  "
$ws.Range("E9").Value = "  +0.37%  "
$ws.Range("E10").Value = "  -1.12%  "
$ws.Range("D11").Value = "5.40"
$ws.Range("E11").Value = "  +0.72%  "
$ws.Range("E12").Value = "  -0.61%  "
$ws.Range("D13").Value = "2.953.16"
$ws.Range("E13").Value = "  +0.46%  "
$ws.Range("D14").Value = "58.898.88"
$ws.Range("E14").Value = "  +0.19%  "
$ws.Range("D15").Value = "22.76"
$ws.Range("E15").Value = "  -1.86%  "
$ws.Range("E16").Value = "  -1.08%  "
$ws.Range("D17").Value = "2.510.70"
$ws.Range("E17").Value = "  +0.11%  "
$ws.Range("D18").Value = "11.06"
$ws.Range("E18").Value = "  +0.08%  "
$ws.Range("D19").Value = "4.25"
$ws.Range("E19").Value = "  -0.07%  "
$ws.Range("D20").Value = "323.86"
$ws.Range("E20").Value = "  -0.19%  "
$ws.Range("E21").Value = "  -0.03%  "
$ws.Range("D22").Value = "5.93"
$ws.Range("E22").Value = "  +1.11%  "
$ws.Range("D23").Value = "64.95"
$ws.Range("E23").Value = "  +0.40%  "
$ws.Range("E24").Value = "  +0.10%  "
$ws.Range("D25").Value = "0.164"
$ws.Range("E25").Value = "  -0.64%  "
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  +0.38%  "
$ws.Range("D27").Value = "7.51"
$ws.Range("E27").Value = "  -0.93%  "
$ws.Range("D28").Value = "0.0₃0762"
$ws.Range("E28").Value = "  -1.20%  "
$ws.Range("D29").Value = "6.46"
$ws.Range("E29").Value = "  -4.53%  "
$ws.Range("E30").Value = "  -1.07%  "
$ws.Range("D31").Value = "168.86"
$ws.Range("E31").Value = "  +0.62%  "
$ws.Range("E33").Value = "  -4.02%  "
$ws.Range("D34").Value = "18.36"
$ws.Range("E34").Value = "  -1.09%  "
$ws.Range("E35").Value = "  -2.81%  "
$ws.Range("D36").Value = "4.04"
$ws.Range("E36").Value = "  -1.41%  "
$ws.Range("E37").Value = "  -3.29%  "
$ws.Range("E38").Value = "  -1.37%  "
$ws.Range("D39").Value = "0.798"
$ws.Range("E39").Value = "  -3.41%  "
$ws.Range("D40").Value = "281.06"
$ws.Range("E40").Value = "  +0.42%  "
$ws.Range("E41").Value = "  +0.32%  "
$ws.Range("D42").Value = "0.603"
$ws.Range("E42").Value = "  -0.31%  "
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").Value = "4.98"
$ws.Range("E43").Value = "  -5.42%  "
$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").Value = "129.61"
$ws.Range("E44").Value = "  +1.11%  "
$ws.Range("B45").Value = "WhiteBITCoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D45").Value = "10.91"
$ws.Range("E45").Value = "  +0.39%  "
$ws.Range("D46").Value = "0.0924"
$ws.Range("E46").Value = "  -0.68%  "
$ws.Range("D47").Value = "0.0500"
$ws.Range("E47").Value = "  -2.87%  "
$ws.Range("E48").Value = "  -1.46%  "
$ws.Range("D49").Value = "17.28"
$ws.Range("E49").Value = "  -0.43%  "
$ws.Range("D50").Value = "1.750.85"
$ws.Range("E50").Value = "  -1.15%  "
$ws.Range("E51").Value = "  -0.54%  "

# Restore default (no explicit number format) style so cells match original plain-text styling
$ws.Range("D2:E51").Style = $ws.Range("B2").Style
